$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '88.505.79'
$ws.Range('E2').Value = '  +0.67%  '

$ws.Range('D3').Value = '3.286.57'
$ws.Range('E3').Value = '  -1.28%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.94'
$ws.Range('E5').Value = '  -2.61%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '629.05'
$ws.Range('E6').Value = '  -1.26%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.379'
$ws.Range('E7').Value = '  +14.81%  '

$ws.Range('E8').Value = '  +17.73%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  +0.06%  '

$ws.Range('D10').Value = '3.279.84'
$ws.Range('E10').Value = '  -1.33%  '

$ws.Range('E11').Value = '  -5.39%  '

$ws.Range('E12').Value = '  +11.94%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000264'
$ws.Range('E13').Value = '  -4.56%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.52'
$ws.Range('E14').Value = '  +2.14%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '34.41'
$ws.Range('E15').Value = '  +0.18%  '

$ws.Range('D16').Value = '3.889.58'
$ws.Range('E16').Value = '  -1.34%  '

$ws.Range('D17').Value = '88.342.56'
$ws.Range('E17').Value = '  +0.93%  '

$ws.Range('D18').Value = '3.300.42'
$ws.Range('E18').Value = '  -0.64%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.19'
$ws.Range('E19').Value = '  -0.70%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.11'
$ws.Range('E20').Value = '  -2.91%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '438.79'
$ws.Range('E21').Value = '  -2.00%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.93'
$ws.Range('E22').Value = '  -1.85%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.37'
$ws.Range('E23').Value = '  +0.48%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.46'
$ws.Range('E24').Value = '  +0.84%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.35'
$ws.Range('E25').Value = '  +0.72%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '5.27'
$ws.Range('E26').Value = '  -1.53%  '

$ws.Range('E27').Value = '  -1.28%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '77.28'
$ws.Range('E28').Value = '  -1.61%  '

$ws.Range('E29').Value = '  +4.86%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.11%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.179'
$ws.Range('E31').Value = '  -3.36%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  -0.20%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '576.16'
$ws.Range('E33').Value = '  -3.93%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '8.81'
$ws.Range('E34').Value = '  -5.48%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.39'
$ws.Range('E35').Value = '  -10.04%  '

$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.16'
$ws.Range('E36').Value = '  +8.72%  '

$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.97'
$ws.Range('E37').Value = '  -4.22%  '

$ws.Range('E38').Value = '  -7.00%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '22.85'
$ws.Range('E39').Value = '  -2.41%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '21.83'
$ws.Range('E40').Value = '  +2.09%  '

$ws.Range('E41').Value = '  +0.17%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.12'
$ws.Range('E42').Value = '  +2.08%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.403'
$ws.Range('E43').Value = '  -3.66%  '

$ws.Range('E44').Value = '  -1.42%  '

$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '153.63'
$ws.Range('E46').Value = '  -2.16%  '

$ws.Range('E47').Value = '  +21.60%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '180.96'
$ws.Range('E48').Value = '  -4.21%  '

$ws.Range('E49').Value = '  -2.81%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.33'
$ws.Range('E50').Value = '  -3.35%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.23'
$ws.Range('E51').Value = '  -0.63%  '
